$d = $word.ActiveDocument

# --- Locate the anchor paragraphs by their current (pre-edit) text ---
# Paragraph "6 aula -> software privado" is the last of the numbered "aula" paragraphs.
# Right after it comes the paragraph that holds only the _GoBack bookmark, and after
# that the single trailing empty paragraph.
$p6 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -match "^6 aula") {
        $p6 = $d.Paragraphs($i)
        break
    }
}
if ($null -eq $p6) {
    throw "Could not locate the '6 aula' paragraph to anchor the new content after."
}

# --- Insert the new "7 aula" paragraph right after it ---
$p6.Range.InsertParagraphAfter() | Out-Null
$pNew1 = $p6.Next()

$frag1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>7</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> aula </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> será começado a aprender o libre office </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>writer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> o que levará em torno de 5 aulas</w:t></w:r></w:p>'
$pNew1.Range.InsertXML($frag1) | Out-Null

# --- Insert the "Após será aprendido..." paragraph right after that one ---
$pNew1.Range.InsertParagraphAfter() | Out-Null
$pNew2 = $pNew1.Next()

$frag2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Após será aprendido o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Impress</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que é para fazer apresentação de slides</w:t></w:r></w:p>'
$pNew2.Range.InsertXML($frag2) | Out-Null

# --- The paragraph that used to contain only the _GoBack bookmark now gets the
#     "A planilha calc..." text woven around the (preserved) bookmark, plus a
#     trailing space run after the bookmark. ---
$pBookmark = $pNew2.Next()

$frag3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">A planilha </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>calc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> quando </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>ser</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ensinada levara em torno de 10 aulas até concluirmos o assunto da mesma pois é muito extenso e tem fórmulas importantes para o aprendizado dos alunos.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$pBookmark.Range.InsertXML($frag3) | Out-Null

# --- Finally, add one new empty paragraph right after the bookmark paragraph
#     (the document already ends with a single empty paragraph; now there are two). ---
$pBookmark.Range.InsertParagraphAfter() | Out-Null
$pEmpty = $pBookmark.Next()
$fragEmpty = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$pEmpty.Range.InsertXML($fragEmpty) | Out-Null

Write-Output "done"
